$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.080.04'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '3.504.98'
$ws.Range("E3").Value = '  +1.31%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = "'572.72"

$ws.Range("D6").Value = "'178.34"
$ws.Range("E6").Value = '  -3.93%  '

$ws.Range("E7").Value = '  +6.09%  '

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = "'0.633"
$ws.Range("E9").Value = '  +1.64%  '

$ws.Range("E10").Value = '  +5.65%  '

$ws.Range("D11").Value = "'55.43"
$ws.Range("E11").Value = '  +2.59%  '

$ws.Range("E12").Value = '  +3.79%  '

$ws.Range("D13").Value = "'9.26"
$ws.Range("E13").Value = '  +0.27%  '

$ws.Range("D14").Value = '4.067.32'
$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("D15").Value = '3.500.67'
$ws.Range("E15").Value = '  +0.76%  '

$ws.Range("E16").Value = '  +0.53%  '

$ws.Range("D17").Value = "'18.35"
$ws.Range("E17").Value = '  +2.26%  '

$ws.Range("D18").Value = '66.054.84'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").Value = "'12.02"
$ws.Range("E19").Value = '  +3.60%  '

$ws.Range("D20").Value = "'1.00"
$ws.Range("E20").Value = '  +2.40%  '

$ws.Range("D21").Value = "'414.04"
$ws.Range("E21").Value = '  +0.97%  '

$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = '  +8.54%  '

$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").Value = "'4.31"
$ws.Range("E23").Value = '  +3.77%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'85.72"
$ws.Range("E24").Value = '  +2.06%  '

$ws.Range("D25").Value = "'13.14"
$ws.Range("E25").Value = '  +12.60%  '

$ws.Range("D26").Value = "'10.98"
$ws.Range("E26").Value = '  -0.48%  '

$ws.Range("D27").Value = "'2.86"
$ws.Range("E27").Value = '  -0.14%  '

$ws.Range("D28").Value = "'9.10"
$ws.Range("E28").Value = '  +5.41%  '

$ws.Range("D29").Value = "'30.44"
$ws.Range("E29").Value = '  +2.41%  '

$ws.Range("D30").Value = "'627.18"
$ws.Range("E30").Value = '  -3.71%  '

$ws.Range("D31").Value = "'6.53"
$ws.Range("E31").Value = '  -0.57%  '

$ws.Range("D32").Value = "'11.68"
$ws.Range("E32").Value = '  +1.05%  '

$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = '  +1.25%  '

$ws.Range("E34").Value = '  +14.67%  '

$ws.Range("D35").Value = "'59.46"
$ws.Range("E35").Value = '  +1.47%  '

$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0794'
$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").Value = "'37.13"
$ws.Range("E38").Value = '  -2.08%  '

$ws.Range("D39").Value = "'0.380"
$ws.Range("E39").Value = '  -1.05%  '

$ws.Range("D40").Value = '3.257.49'
$ws.Range("E40").Value = '  +9.16%  '

$ws.Range("D41").Value = "'3.41"
$ws.Range("E41").Value = '  +2.82%  '

$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("D43").Value = "'2.92"
$ws.Range("E43").Value = '  +2.36%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = "'0.0418"
$ws.Range("E44").Value = '  +2.08%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = "'3.28"
$ws.Range("E45").Value = '  -3.03%  '

$ws.Range("E46").Value = '  -3.03%  '

$ws.Range("D47").Value = "'2.71"
$ws.Range("E47").Value = '  +0.83%  '

$ws.Range("D48").Value = "'0.132"
$ws.Range("E48").Value = '  +3.27%  '

$ws.Range("D49").Value = "'8.64"
$ws.Range("E49").Value = '  -0.96%  '

$ws.Range("D50").Value = "'140.07"
$ws.Range("E50").Value = '  +0.78%  '

$ws.Range("D51").Value = "'2.37"
$ws.Range("E51").Value = '  +0.68%  '
